$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark off its old location (it will be
#    re-added at the new edit point below, mirroring how Word tracks
#    the most-recent-edit bookmark).
# ------------------------------------------------------------------
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# ------------------------------------------------------------------
# 2. Turn the original first bullet ("Install a Ruby IDE or text
#    editor:") into the new "Install Git:" bullet, then insert two
#    fresh bullets after it:
#       - a sub-bullet linking to "Set Up Git"
#       - a restored "Install a Ruby IDE or text editor:" bullet
# ------------------------------------------------------------------
$p = $d.Paragraphs(2)

# Use a temporary suffix so the zero-width bookmark we add lands on an
# ordinary mid-run boundary rather than the exact "end of paragraph"
# position (collapsing a range there and handing it to Bookmarks.Add
# mis-resolves in this host, so we sidestep it: add the bookmark next
# to temporary filler text, then delete the filler so the bookmark
# collapses to the correct spot right before the paragraph mark).
$p.Range.Text = "Install Git:ZZZ"
$p2 = $d.Paragraphs(2)
$cut = $p2.Range.Start + 12
$bmSpot = $d.Range($cut, $cut)
$d.Bookmarks.Add("_GoBack", $bmSpot)
$filler = $d.Range($cut, $cut + 3)
$filler.Delete()

# New sub-bullet: "Set Up Git" hyperlink (ilvl 1)
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$gitPara = $d.Paragraphs(3)
$gitPara.Range.Text = "Set Up Git"
$gitPara = $d.Paragraphs(3)
$linkRange = $d.Range($gitPara.Range.Start, $gitPara.Range.Start + 10)
$d.Hyperlinks.Add($linkRange, "https://help.github.com/articles/set-up-git/")
$gitPara = $d.Paragraphs(3)
$gitPara.Range.ListFormat.ListIndent()

# New bullet restoring "Install a Ruby IDE or text editor:" (ilvl 0)
$gitPara = $d.Paragraphs(3)
$gitPara.Range.InsertParagraphAfter()
$idePara = $d.Paragraphs(4)
$idePara.Range.Text = "Install a Ruby IDE or text editor:"
$idePara = $d.Paragraphs(4)
$idePara.Range.ListFormat.ListOutdent()

Write-Output "Para2: $($d.Paragraphs(2).Range.Text)"
Write-Output "Para3: $($d.Paragraphs(3).Range.Text)"
Write-Output "Para4: $($d.Paragraphs(4).Range.Text)"
Write-Output "Para5: $($d.Paragraphs(5).Range.Text)"
